$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cells I1 and J1 (same style as existing headers, e.g. style index 1)
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)

# I and J data values, row 2 through 73
$ijValues = @{
    2 = @(7, 8)
    3 = @(8, 8)
    4 = @(8, 8)
    5 = @(1, 3)
    6 = @(6, 6)
    7 = @(6, 6)
    8 = @(6, 7)
    9 = @(7, 7)
    10 = @(7, 8)
    11 = @(6, 7)
    12 = @(7, 7)
    13 = @(6, 7)
    14 = @(7, 7)
    15 = @(6, 6)
    16 = @(7, 7)
    17 = @(6, 6)
    18 = @(1, 3)
    19 = @(8, 8)
    20 = @(9, 9)
    21 = @(8, 8)
    22 = @(6, 6)
    23 = @(9, 10)
    24 = @(6, 6)
    25 = @(7, 7)
    26 = @(6, 6)
    27 = @(9, 10)
    28 = @(5, 6)
    29 = @(6, 6)
    30 = @(7, 7)
    31 = @(8, 8)
    32 = @(5, 6)
    33 = @(7, 7)
    34 = @(9, 9)
    35 = @(7, 7)
    36 = @(7, 8)
    37 = @(9, 9)
    38 = @(4, 5)
    39 = @(10, 10)
    40 = @(8, 8)
    41 = @(6, 6)
    42 = @(5, 6)
    43 = @(8, 9)
    44 = @(6, 6)
    45 = @(4, 5)
    46 = @(9, 9)
    47 = @(2, 2)
    48 = @(8, 8)
    49 = @(6, 6)
    50 = @(6, 6)
    51 = @(6, 7)
    52 = @(8, 8)
    53 = @(6, 6)
    54 = @(6, 6)
    55 = @(8, 8)
    56 = @(4, 5)
    57 = @(8, 8)
    58 = @(5, 6)
    59 = @(8, 8)
    60 = @(8, 8)
    61 = @(6, 7)
    62 = @(7, 7)
    63 = @(8, 8)
    64 = @(6, 6)
    65 = @(5, 5)
    66 = @(7, 7)
    67 = @(9, 9)
    68 = @(6, 7)
    69 = @(6, 7)
    70 = @(6, 6)
    71 = @(8, 8)
    72 = @(5, 6)
    73 = @(3, 3)
}

foreach ($row in $ijValues.Keys) {
    $vals = $ijValues[$row]
    $ws.Cells.Item($row, 9).Value = $vals[0]
    $ws.Cells.Item($row, 10).Value = $vals[1]
}
